$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (I2:L2) - M1_PH, CM2_PH, CMN3_PH, CMN4_PH
$ws.Range("I2").Value = -0.2051071372574013
$ws.Range("J2").Value = 0.3037700709882454
$ws.Range("K2").Value = -0.1987230108742328
$ws.Range("L2").Value = 2.292162609111811

# Row 19 (I19:L19) - M1_PH, CM2_PH, CMN3_PH, CMN4_PH
$ws.Range("I19").Value = -0.4335568406464395
$ws.Range("J19").Value = 0.397851208992047
$ws.Range("K19").Value = 0.1273926462583985
$ws.Range("L19").Value = 1.916091474589909
